$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 263.5  # H12
$ws.Cells.Item(12, 9).Value = 195  # I12
$ws.Cells.Item(12, 10).Value = 400.5  # J12
$ws.Cells.Item(12, 11).Value = 195  # K12
$ws.Cells.Item(12, 12).Value = 400.5  # L12
$ws.Cells.Item(12, 13).Value = -25  # M12
$ws.Cells.Item(12, 14).Value = -740.5  # N12
$ws.Cells.Item(15, 8).Value = 1584.0857  # H15
$ws.Cells.Item(15, 9).Value = 1584.0857  # I15
$ws.Cells.Item(15, 11).Value = 4752.257100000001  # K15
$ws.Cells.Item(15, 13).Value = -4583.257100000001  # M15
$ws.Cells.Item(17, 8).Value = 12776.667  # H17
$ws.Cells.Item(17, 9).Value = 1500  # I17
$ws.Cells.Item(17, 10).Value = 14186.25  # J17
$ws.Cells.Item(17, 11).Value = 4500  # K17
$ws.Cells.Item(17, 12).Value = 42558.75  # L17
$ws.Cells.Item(17, 13).Value = -4332  # M17
$ws.Cells.Item(17, 14).Value = -42894.75  # N17
$ws.Cells.Item(28, 8).Value = 1053.5714  # H28
$ws.Cells.Item(28, 9).Value = 364.92307  # I28
$ws.Cells.Item(28, 11).Value = 364.92307  # K28
$ws.Cells.Item(28, 13).Value = 120.07693  # M28
$ws.Cells.Item(39, 8).Value = 231.55556  # H39
$ws.Cells.Item(39, 9).Value = 231.55556  # I39
$ws.Cells.Item(39, 11).Value = 694.66668  # K39
$ws.Cells.Item(39, 13).Value = -398.66668  # M39
$ws.Cells.Item(42, 8).Value = 410.85715  # H42
$ws.Cells.Item(42, 9).Value = 314.83334  # I42
$ws.Cells.Item(42, 11).Value = 944.5000200000001  # K42
$ws.Cells.Item(42, 13).Value = -714.5000200000001  # M42
$ws.Cells.Item(70, 8).Value = 22338610  # H70
$ws.Cells.Item(70, 9).Value = 252624.75  # I70
$ws.Cells.Item(70, 10).Value = 40007400  # J70
$ws.Cells.Item(70, 11).Value = 757874.25  # K70
$ws.Cells.Item(70, 12).Value = 120022200  # L70
$ws.Cells.Item(70, 13).Value = -757604.25  # M70
$ws.Cells.Item(70, 14).Value = -120022740  # N70
$ws.Cells.Item(73, 8).Value = 22338610  # H73
$ws.Cells.Item(73, 9).Value = 252624.75  # I73
$ws.Cells.Item(73, 10).Value = 40007400  # J73
$ws.Cells.Item(73, 11).Value = 757874.25  # K73
$ws.Cells.Item(73, 12).Value = 120022200  # L73
$ws.Cells.Item(73, 13).Value = -756938.25  # M73
$ws.Cells.Item(73, 14).Value = -120024072  # N73
$ws.Cells.Item(76, 8).Value = 7296.731  # H76
$ws.Cells.Item(76, 9).Value = 6350.6875  # I76
$ws.Cells.Item(76, 10).Value = 8810.4  # J76
$ws.Cells.Item(76, 11).Value = 6350.6875  # K76
$ws.Cells.Item(76, 12).Value = 8810.4  # L76
$ws.Cells.Item(76, 13).Value = -6035.6875  # M76
$ws.Cells.Item(76, 14).Value = -9440.4  # N76
$ws.Cells.Item(79, 8).Value = 7296.731  # H79
$ws.Cells.Item(79, 9).Value = 6350.6875  # I79
$ws.Cells.Item(79, 10).Value = 8810.4  # J79
$ws.Cells.Item(79, 11).Value = 6350.6875  # K79
$ws.Cells.Item(79, 12).Value = 8810.4  # L79
$ws.Cells.Item(79, 13).Value = -5258.6875  # M79
$ws.Cells.Item(79, 14).Value = -10994.4  # N79
$ws.Cells.Item(96, 8).Value = 2063.4285  # H96
$ws.Cells.Item(96, 9).Value = 2983  # I96
$ws.Cells.Item(96, 10).Value = 1373.75  # J96
$ws.Cells.Item(96, 11).Value = 8949  # K96
$ws.Cells.Item(96, 12).Value = 4121.25  # L96
$ws.Cells.Item(96, 13).Value = -7576  # M96
$ws.Cells.Item(96, 14).Value = -6867.25  # N96
$ws.Cells.Item(98, 8).Value = 437899.22  # H98
$ws.Cells.Item(98, 9).Value = 2235.2727  # I98
$ws.Cells.Item(98, 10).Value = 2035333.6  # J98
$ws.Cells.Item(98, 11).Value = 2235.2727  # K98
$ws.Cells.Item(98, 12).Value = 2035333.6  # L98
$ws.Cells.Item(98, 13).Value = -737.2727  # M98
$ws.Cells.Item(98, 14).Value = -2038329.6  # N98
$ws.Cells.Item(100, 8).Value = 4726.636  # H100
$ws.Cells.Item(100, 10).Value = 6191.8335  # J100
$ws.Cells.Item(100, 12).Value = 6191.8335  # L100
$ws.Cells.Item(100, 14).Value = -7273.8335  # N100
$ws.Cells.Item(103, 8).Value = 1479.7222  # H103
$ws.Cells.Item(103, 9).Value = 360.77777  # I103
$ws.Cells.Item(103, 10).Value = 2598.6667  # J103
$ws.Cells.Item(103, 11).Value = 1082.33331  # K103
$ws.Cells.Item(103, 12).Value = 7796.000100000001  # L103
$ws.Cells.Item(103, 13).Value = -496.33331  # M103
$ws.Cells.Item(103, 14).Value = -8968.000100000001  # N103
$ws.Cells.Item(112, 8).Value = 2020.591  # H112
$ws.Cells.Item(112, 10).Value = 2048.9524  # J112
$ws.Cells.Item(112, 12).Value = 6146.8572  # L112
$ws.Cells.Item(112, 14).Value = -8362.8572  # N112
$ws.Cells.Item(113, 8).Value = 7371.143  # H113
$ws.Cells.Item(113, 9).Value = 4308  # I113
$ws.Cells.Item(113, 11).Value = 4308  # K113
$ws.Cells.Item(113, 13).Value = -1054  # M113
$ws.Cells.Item(122, 8).Value = 437899.22  # H122
$ws.Cells.Item(122, 9).Value = 2235.2727  # I122
$ws.Cells.Item(122, 10).Value = 2035333.6  # J122
$ws.Cells.Item(122, 11).Value = 6705.8181  # K122
$ws.Cells.Item(122, 12).Value = 6106000.800000001  # L122
$ws.Cells.Item(122, 13).Value = -4255.8181  # M122
$ws.Cells.Item(122, 14).Value = -6110900.800000001  # N122
$ws.Cells.Item(125, 8).Value = 2614  # H125
$ws.Cells.Item(125, 9).Value = 900  # I125
$ws.Cells.Item(125, 10).Value = 3042.5  # J125
$ws.Cells.Item(125, 11).Value = 8100  # K125
$ws.Cells.Item(125, 12).Value = 27382.5  # L125
$ws.Cells.Item(125, 13).Value = -5640  # M125
$ws.Cells.Item(125, 14).Value = -32302.5  # N125
$ws.Cells.Item(138, 8).Value = 3367  # H138
$ws.Cells.Item(138, 9).Value = 2223.889  # I138
$ws.Cells.Item(138, 10).Value = 3786.9185  # J138
$ws.Cells.Item(138, 11).Value = 6671.667  # K138
$ws.Cells.Item(138, 12).Value = 11360.7555  # L138
$ws.Cells.Item(138, 13).Value = -1531.667  # M138
$ws.Cells.Item(138, 14).Value = -21640.7555  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 4340.5483  # H2
$ws.Cells.Item(2, 9).Value = 781.03705  # I2
$ws.Cells.Item(2, 10).Value = 28367.25  # J2
$ws.Cells.Item(2, 11).Value = 781.03705  # K2
$ws.Cells.Item(2, 12).Value = 28367.25  # L2
$ws.Cells.Item(2, 13).Value = -668.03705  # M2
$ws.Cells.Item(2, 14).Value = -28593.25  # N2
$ws.Cells.Item(22, 8).Value = 7750  # H22
$ws.Cells.Item(22, 10).Value = 7750  # J22
$ws.Cells.Item(22, 12).Value = 7750  # L22
$ws.Cells.Item(22, 14).Value = -8348  # N22
$ws.Cells.Item(32, 8).Value = 4762.7144  # H32
$ws.Cells.Item(32, 9).Value = 4475.625  # I32
$ws.Cells.Item(32, 10).Value = 10504.5  # J32
$ws.Cells.Item(32, 11).Value = 4475.625  # K32
$ws.Cells.Item(32, 12).Value = 10504.5  # L32
$ws.Cells.Item(32, 13).Value = -4188.625  # M32
$ws.Cells.Item(32, 14).Value = -11078.5  # N32
$ws.Cells.Item(45, 8).Value = 3240.2942  # H45
$ws.Cells.Item(45, 9).Value = 2590.0833  # I45
$ws.Cells.Item(45, 11).Value = 2590.0833  # K45
$ws.Cells.Item(45, 13).Value = -2213.0833  # M45
$ws.Cells.Item(61, 8).Value = 3418.7104  # H61
$ws.Cells.Item(61, 9).Value = 2796.3704  # I61
$ws.Cells.Item(61, 11).Value = 2796.3704  # K61
$ws.Cells.Item(61, 13).Value = -2584.3704  # M61
$ws.Cells.Item(63, 8).Value = 6280.5625  # H63
$ws.Cells.Item(63, 9).Value = 3061.75  # I63
$ws.Cells.Item(63, 10).Value = 9499.375  # J63
$ws.Cells.Item(63, 11).Value = 3061.75  # K63
$ws.Cells.Item(63, 12).Value = 9499.375  # L63
$ws.Cells.Item(63, 13).Value = -2375.75  # M63
$ws.Cells.Item(63, 14).Value = -10871.375  # N63
$ws.Cells.Item(66, 8).Value = 6280.5625  # H66
$ws.Cells.Item(66, 9).Value = 3061.75  # I66
$ws.Cells.Item(66, 10).Value = 9499.375  # J66
$ws.Cells.Item(66, 11).Value = 15308.75  # K66
$ws.Cells.Item(66, 12).Value = 47496.875  # L66
$ws.Cells.Item(66, 13).Value = -11876.75  # M66
$ws.Cells.Item(66, 14).Value = -54360.875  # N66
$ws.Cells.Item(116, 8).Value = 4340.5483  # H116
$ws.Cells.Item(116, 9).Value = 781.03705  # I116
$ws.Cells.Item(116, 10).Value = 28367.25  # J116
$ws.Cells.Item(116, 11).Value = 781.03705  # K116
$ws.Cells.Item(116, 12).Value = 28367.25  # L116
$ws.Cells.Item(116, 13).Value = 1512.96295  # M116
$ws.Cells.Item(116, 14).Value = -32955.25  # N116
$ws.Cells.Item(122, 8).Value = 4573.2666  # H122
$ws.Cells.Item(122, 9).Value = 4585.7144  # I122
$ws.Cells.Item(122, 10).Value = 4562.375  # J122
$ws.Cells.Item(122, 11).Value = 13757.1432  # K122
$ws.Cells.Item(122, 12).Value = 13687.125  # L122
$ws.Cells.Item(122, 13).Value = -11307.1432  # M122
$ws.Cells.Item(122, 14).Value = -18587.125  # N122
$ws.Cells.Item(132, 8).Value = 3000.3333  # H132
$ws.Cells.Item(132, 9).Value = 2464.0476  # I132
$ws.Cells.Item(132, 10).Value = 3938.8333  # J132
$ws.Cells.Item(132, 11).Value = 7392.1428  # K132
$ws.Cells.Item(132, 12).Value = 11816.4999  # L132
$ws.Cells.Item(132, 13).Value = -4862.1428  # M132
$ws.Cells.Item(132, 14).Value = -16876.4999  # N132
$ws.Cells.Item(136, 8).Value = 3418.7104  # H136
$ws.Cells.Item(136, 9).Value = 2796.3704  # I136
$ws.Cells.Item(136, 11).Value = 8389.111199999999  # K136
$ws.Cells.Item(136, 13).Value = -5839.111199999999  # M136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 4340.5483  # H3
$ws.Cells.Item(3, 9).Value = 781.03705  # I3
$ws.Cells.Item(3, 10).Value = 28367.25  # J3
$ws.Cells.Item(3, 11).Value = 781.03705  # K3
$ws.Cells.Item(3, 12).Value = 28367.25  # L3
$ws.Cells.Item(3, 13).Value = -667.03705  # M3
$ws.Cells.Item(3, 14).Value = -28595.25  # N3
$ws.Cells.Item(86, 8).Value = 5482.923  # H86
$ws.Cells.Item(86, 9).Value = 4638.636  # I86
$ws.Cells.Item(86, 11).Value = 4638.636  # K86
$ws.Cells.Item(86, 13).Value = -3515.636  # M86
$ws.Cells.Item(89, 8).Value = 5482.923  # H89
$ws.Cells.Item(89, 9).Value = 4638.636  # I89
$ws.Cells.Item(89, 11).Value = 23193.18  # K89
$ws.Cells.Item(89, 13).Value = -17577.18  # M89
$ws.Cells.Item(105, 8).Value = 15256.92  # H105
$ws.Cells.Item(105, 9).Value = 14562.3125  # I105
$ws.Cells.Item(105, 10).Value = 16491.777  # J105
$ws.Cells.Item(105, 11).Value = 14562.3125  # K105
$ws.Cells.Item(105, 12).Value = 16491.777  # L105
$ws.Cells.Item(105, 13).Value = -12815.3125  # M105
$ws.Cells.Item(105, 14).Value = -19985.777  # N105

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2421.4  # H16
$ws.Cells.Item(16, 9).Value = 1977.1  # I16
$ws.Cells.Item(16, 10).Value = 3310  # J16
$ws.Cells.Item(16, 11).Value = 1977.1  # K16
$ws.Cells.Item(16, 12).Value = 3310  # L16
$ws.Cells.Item(16, 13).Value = -1690.1  # M16
$ws.Cells.Item(16, 14).Value = -3884  # N16
$ws.Cells.Item(28, 8).Value = 10643  # H28
$ws.Cells.Item(28, 10).Value = 10643  # J28
$ws.Cells.Item(28, 12).Value = 10643  # L28
$ws.Cells.Item(28, 14).Value = -11133  # N28
$ws.Cells.Item(31, 8).Value = 23821.469  # H31
$ws.Cells.Item(31, 9).Value = 2241.4736  # I31
$ws.Cells.Item(31, 11).Value = 2241.4736  # K31
$ws.Cells.Item(31, 13).Value = -1946.4736  # M31
$ws.Cells.Item(34, 8).Value = 23821.469  # H34
$ws.Cells.Item(34, 9).Value = 2241.4736  # I34
$ws.Cells.Item(34, 11).Value = 2241.4736  # K34
$ws.Cells.Item(34, 13).Value = -2039.4736  # M34
$ws.Cells.Item(43, 8).Value = 39000  # H43
$ws.Cells.Item(43, 10).Value = 39000  # J43
$ws.Cells.Item(43, 12).Value = 39000  # L43
$ws.Cells.Item(43, 14).Value = -39368  # N43
$ws.Cells.Item(56, 8).Value = 40000  # H56
$ws.Cells.Item(56, 9).Value = 20000  # I56
$ws.Cells.Item(56, 11).Value = 20000  # K56
$ws.Cells.Item(56, 13).Value = -19155  # M56
$ws.Cells.Item(62, 8).Value = 9348.666999999999  # H62
$ws.Cells.Item(62, 10).Value = 29003  # J62
$ws.Cells.Item(62, 12).Value = 29003  # L62
$ws.Cells.Item(62, 14).Value = -30251  # N62
$ws.Cells.Item(65, 8).Value = 9348.666999999999  # H65
$ws.Cells.Item(65, 10).Value = 29003  # J65
$ws.Cells.Item(65, 12).Value = 145015  # L65
$ws.Cells.Item(65, 14).Value = -151255  # N65
$ws.Cells.Item(68, 8).Value = 50000  # H68
$ws.Cells.Item(68, 10).Value = 50000  # J68
$ws.Cells.Item(68, 12).Value = 50000  # L68
$ws.Cells.Item(68, 14).Value = -51498  # N68
$ws.Cells.Item(71, 8).Value = 50000  # H71
$ws.Cells.Item(71, 10).Value = 50000  # J71
$ws.Cells.Item(71, 12).Value = 150000  # L71
$ws.Cells.Item(71, 14).Value = -157488  # N71
$ws.Cells.Item(94, 8).Value = 5642.143  # H94
$ws.Cells.Item(94, 9).Value = 4376.75  # I94
$ws.Cells.Item(94, 10).Value = 7329.3335  # J94
$ws.Cells.Item(94, 11).Value = 4376.75  # K94
$ws.Cells.Item(94, 12).Value = 7329.3335  # L94
$ws.Cells.Item(94, 13).Value = -3925.75  # M94
$ws.Cells.Item(94, 14).Value = -8231.333500000001  # N94
$ws.Cells.Item(101, 8).Value = 39000  # H101
$ws.Cells.Item(101, 10).Value = 39000  # J101
$ws.Cells.Item(101, 12).Value = 39000  # L101
$ws.Cells.Item(101, 14).Value = -45490  # N101
$ws.Cells.Item(113, 8).Value = 2421.4  # H113
$ws.Cells.Item(113, 9).Value = 1977.1  # I113
$ws.Cells.Item(113, 10).Value = 3310  # J113
$ws.Cells.Item(113, 11).Value = 1977.1  # K113
$ws.Cells.Item(113, 12).Value = 3310  # L113
$ws.Cells.Item(113, 13).Value = 192.9000000000001  # M113
$ws.Cells.Item(113, 14).Value = -7650  # N113
$ws.Cells.Item(122, 8).Value = 3754.52  # H122
$ws.Cells.Item(122, 9).Value = 1019.1579  # I122
$ws.Cells.Item(122, 11).Value = 3057.4737  # K122
$ws.Cells.Item(122, 13).Value = -607.4737  # M122

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(63, 8).Value = 10000  # H63
$ws.Cells.Item(63, 9).Value = 10000  # I63
$ws.Cells.Item(63, 11).Value = 30000  # K63
$ws.Cells.Item(63, 13).Value = -29251  # M63
$ws.Cells.Item(66, 8).Value = 10000  # H66
$ws.Cells.Item(66, 9).Value = 10000  # I66
$ws.Cells.Item(66, 11).Value = 90000  # K66
$ws.Cells.Item(66, 13).Value = -86256  # M66
$ws.Cells.Item(68, 8).Value = 1908.1818  # H68
$ws.Cells.Item(68, 9).Value = 1447.8  # I68
$ws.Cells.Item(68, 10).Value = 2291.8333  # J68
$ws.Cells.Item(68, 11).Value = 4343.4  # K68
$ws.Cells.Item(68, 12).Value = 6875.499899999999  # L68
$ws.Cells.Item(68, 13).Value = -3532.4  # M68
$ws.Cells.Item(68, 14).Value = -8497.499899999999  # N68
$ws.Cells.Item(71, 8).Value = 1908.1818  # H71
$ws.Cells.Item(71, 9).Value = 1447.8  # I71
$ws.Cells.Item(71, 10).Value = 2291.8333  # J71
$ws.Cells.Item(71, 11).Value = 13030.2  # K71
$ws.Cells.Item(71, 12).Value = 20626.4997  # L71
$ws.Cells.Item(71, 13).Value = -8974.199999999999  # M71
$ws.Cells.Item(71, 14).Value = -28738.4997  # N71
$ws.Cells.Item(92, 8).Value = 500  # H92
$ws.Cells.Item(92, 9).Value = 500  # I92
$ws.Cells.Item(92, 10).Value = 0  # J92
$ws.Cells.Item(92, 11).Value = 1500  # K92
$ws.Cells.Item(92, 12).Value = 0  # L92
$ws.Cells.Item(92, 13).Value = -252  # M92
$ws.Cells.Item(92, 14).Value = $null  # N92
$ws.Cells.Item(107, 8).Value = 261.5  # H107
$ws.Cells.Item(107, 9).Value = 208.78572  # I107
$ws.Cells.Item(107, 11).Value = 626.35716  # K107
$ws.Cells.Item(107, 13).Value = 1293.64284  # M107
$ws.Cells.Item(132, 8).Value = 3755.4075  # H132
$ws.Cells.Item(132, 10).Value = 4404.4614  # J132
$ws.Cells.Item(132, 12).Value = 39640.1526  # L132
$ws.Cells.Item(132, 14).Value = -44700.1526  # N132

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(6, 8).Value = 3032.25  # H6
$ws.Cells.Item(6, 9).Value = 739.5  # I6
$ws.Cells.Item(6, 10).Value = 5325  # J6
$ws.Cells.Item(6, 11).Value = 739.5  # K6
$ws.Cells.Item(6, 12).Value = 5325  # L6
$ws.Cells.Item(6, 13).Value = -626.5  # M6
$ws.Cells.Item(6, 14).Value = -5551  # N6
$ws.Cells.Item(11, 8).Value = 7134438  # H11
$ws.Cells.Item(11, 9).Value = 11666666  # I11
$ws.Cells.Item(11, 11).Value = 11666666  # K11
$ws.Cells.Item(11, 13).Value = -11666527  # M11
$ws.Cells.Item(16, 8).Value = 3032.25  # H16
$ws.Cells.Item(16, 9).Value = 739.5  # I16
$ws.Cells.Item(16, 10).Value = 5325  # J16
$ws.Cells.Item(16, 11).Value = 739.5  # K16
$ws.Cells.Item(16, 12).Value = 5325  # L16
$ws.Cells.Item(16, 13).Value = -489.5  # M16
$ws.Cells.Item(16, 14).Value = -5825  # N16
$ws.Cells.Item(57, 8).Value = 45000  # H57
$ws.Cells.Item(57, 10).Value = 45000  # J57
$ws.Cells.Item(57, 12).Value = 45000  # L57
$ws.Cells.Item(57, 14).Value = -46640  # N57
$ws.Cells.Item(80, 8).Value = 212050.25  # H80
$ws.Cells.Item(80, 9).Value = 386753.62  # I80
$ws.Cells.Item(80, 10).Value = 5582.636  # J80
$ws.Cells.Item(80, 11).Value = 386753.62  # K80
$ws.Cells.Item(80, 12).Value = 5582.636  # L80
$ws.Cells.Item(80, 13).Value = -385755.62  # M80
$ws.Cells.Item(80, 14).Value = -7578.636  # N80
$ws.Cells.Item(83, 8).Value = 212050.25  # H83
$ws.Cells.Item(83, 9).Value = 386753.62  # I83
$ws.Cells.Item(83, 10).Value = 5582.636  # J83
$ws.Cells.Item(83, 11).Value = 1933768.1  # K83
$ws.Cells.Item(83, 12).Value = 27913.18  # L83
$ws.Cells.Item(83, 13).Value = -1928776.1  # M83
$ws.Cells.Item(83, 14).Value = -37897.18  # N83
$ws.Cells.Item(93, 8).Value = 34069.668  # H93
$ws.Cells.Item(93, 9).Value = 28228  # I93
$ws.Cells.Item(93, 10).Value = 34600.727  # J93
$ws.Cells.Item(93, 11).Value = 28228  # K93
$ws.Cells.Item(93, 12).Value = 34600.727  # L93
$ws.Cells.Item(93, 13).Value = -26356  # M93
$ws.Cells.Item(93, 14).Value = -38344.727  # N93
$ws.Cells.Item(102, 8).Value = 2787.5625  # H102
$ws.Cells.Item(102, 9).Value = 1705.8667  # I102
$ws.Cells.Item(102, 11).Value = 1705.8667  # K102
$ws.Cells.Item(102, 13).Value = -83.86670000000004  # M102
$ws.Cells.Item(113, 8).Value = 2971.6365  # H113
$ws.Cells.Item(113, 9).Value = 2429.75  # I113
$ws.Cells.Item(113, 10).Value = 4416.6665  # J113
$ws.Cells.Item(113, 11).Value = 2429.75  # K113
$ws.Cells.Item(113, 12).Value = 4416.6665  # L113
$ws.Cells.Item(113, 13).Value = -259.75  # M113
$ws.Cells.Item(113, 14).Value = -8756.666499999999  # N113
$ws.Cells.Item(122, 8).Value = 24899.375  # H122
$ws.Cells.Item(122, 9).Value = 35499.5  # I122
$ws.Cells.Item(122, 10).Value = 14299.25  # J122
$ws.Cells.Item(122, 11).Value = 106498.5  # K122
$ws.Cells.Item(122, 12).Value = 42897.75  # L122
$ws.Cells.Item(122, 13).Value = -104048.5  # M122
$ws.Cells.Item(122, 14).Value = -47797.75  # N122
$ws.Cells.Item(126, 8).Value = 3974.0527  # H126
$ws.Cells.Item(126, 9).Value = 2233  # I126
$ws.Cells.Item(126, 10).Value = 5541  # J126
$ws.Cells.Item(126, 11).Value = 6699  # K126
$ws.Cells.Item(126, 12).Value = 16623  # L126
$ws.Cells.Item(126, 13).Value = -4229  # M126
$ws.Cells.Item(126, 14).Value = -21563  # N126
$ws.Cells.Item(132, 8).Value = 2939.8635  # H132
$ws.Cells.Item(132, 9).Value = 2581.8857  # I132
$ws.Cells.Item(132, 10).Value = 4332  # J132
$ws.Cells.Item(132, 11).Value = 7745.657099999999  # K132
$ws.Cells.Item(132, 12).Value = 12996  # L132
$ws.Cells.Item(132, 13).Value = -5215.657099999999  # M132
$ws.Cells.Item(132, 14).Value = -18056  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(24, 8).Value = 12650  # H24
$ws.Cells.Item(24, 9).Value = 12650  # I24
$ws.Cells.Item(24, 11).Value = 12650  # K24
$ws.Cells.Item(24, 13).Value = -12307  # M24
$ws.Cells.Item(40, 8).Value = 6330.273  # H40
$ws.Cells.Item(40, 9).Value = 5258.0386  # I40
$ws.Cells.Item(40, 10).Value = 10312.857  # J40
$ws.Cells.Item(40, 11).Value = 5258.0386  # K40
$ws.Cells.Item(40, 12).Value = 10312.857  # L40
$ws.Cells.Item(40, 13).Value = -5122.0386  # M40
$ws.Cells.Item(40, 14).Value = -10584.857  # N40
$ws.Cells.Item(46, 8).Value = 8733.1  # H46
$ws.Cells.Item(46, 10).Value = 10000  # J46
$ws.Cells.Item(46, 12).Value = 10000  # L46
$ws.Cells.Item(46, 14).Value = -10376  # N46
$ws.Cells.Item(122, 8).Value = 507232.5  # H122
$ws.Cells.Item(122, 9).Value = 672442.8  # I122
$ws.Cells.Item(122, 11).Value = 2017328.4  # K122
$ws.Cells.Item(122, 13).Value = -2014878.4  # M122

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 26501  # H18
$ws.Cells.Item(18, 9).Value = 24500  # I18
$ws.Cells.Item(18, 11).Value = 24500  # K18
$ws.Cells.Item(18, 13).Value = -24327  # M18
$ws.Cells.Item(39, 8).Value = 50000000  # H39
$ws.Cells.Item(39, 9).Value = 50000000  # I39
$ws.Cells.Item(39, 11).Value = 50000000  # K39
$ws.Cells.Item(39, 13).Value = -49999587  # M39
$ws.Cells.Item(100, 8).Value = 962.3714  # H100
$ws.Cells.Item(100, 9).Value = 779.44446  # I100
$ws.Cells.Item(100, 11).Value = 1558.88892  # K100
$ws.Cells.Item(100, 13).Value = -1017.88892  # M100
$ws.Cells.Item(113, 8).Value = 346.84  # H113
$ws.Cells.Item(113, 9).Value = 339.26315  # I113
$ws.Cells.Item(113, 11).Value = 1017.78945  # K113
$ws.Cells.Item(113, 13).Value = 1152.21055  # M113
$ws.Cells.Item(122, 8).Value = 5073.5454  # H122
$ws.Cells.Item(122, 9).Value = 3972  # I122
$ws.Cells.Item(122, 11).Value = 11916  # K122
$ws.Cells.Item(122, 13).Value = -9466  # M122
$ws.Cells.Item(136, 8).Value = 4562.75  # H136
$ws.Cells.Item(136, 9).Value = 3159.625  # I136
$ws.Cells.Item(136, 10).Value = 6433.5835  # J136
$ws.Cells.Item(136, 11).Value = 9478.875  # K136
$ws.Cells.Item(136, 12).Value = 19300.7505  # L136
$ws.Cells.Item(136, 13).Value = -6928.875  # M136
$ws.Cells.Item(136, 14).Value = -24400.7505  # N136
